$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 = session 20, plan date "2025.01.13" (col C). Mark this session as completed:
# - remove the yellow "not started" highlight from Duration (B7) and Plan (C7)
# - fill in the Actual date (D7) with the same plan date
# - bump Progress (E7) from 0 to 0.9 (90%)
# - record the Subject (F7) as "Naive Bayes"

$ws.Range("B7").Interior.Pattern = -4142   # xlNone - remove highlight
$ws.Range("C7").Interior.Pattern = -4142   # xlNone - remove highlight

$ws.Range("C7").Copy()
$ws.Range("D7").PasteSpecial(-4163)        # xlPasteValues - copy C7's text as a plain value

$ws.Range("E7").Value = 0.9
$ws.Range("F7").Value = "Na" + [char]0x00EF + "ve Bayes"

[void]$ws.Range("E7").Select()
